$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $text) {
    $r = $ws.Range($cellRef)
    $r.NumberFormat = "@"
    $r.Value = $text
    $r.ClearFormats()
}

Set-TextValue "D2" "29.089.15"
Set-TextValue "E2" "  -1.85%  "

Set-TextValue "D3" "1.834.42"
Set-TextValue "E3" "  -1.42%  "

Set-TextValue "D4" "0.9994"
Set-TextValue "E4" "  -0.02%  "

Set-TextValue "D5" "239.70"
Set-TextValue "E5" "  -2.32%  "

Set-TextValue "D6" "0.6799"
Set-TextValue "E6" "  -2.53%  "

Set-TextValue "D8" "0.2981"
Set-TextValue "E8" "  -2.80%  "

Set-TextValue "D9" "0.07438"
Set-TextValue "E9" "  -3.72%  "

Set-TextValue "D10" "23.15"
Set-TextValue "E10" "  -2.19%  "

Set-TextValue "D11" "0.07647"
Set-TextValue "E11" "  -1.52%  "

Set-TextValue "D12" "1.836.17"
Set-TextValue "E12" "  -1.31%  "

Set-TextValue "D13" "5.018"

Set-TextValue "D14" "0.6756"
Set-TextValue "E14" "  -2.54%  "

Set-TextValue "D15" "86.55"
Set-TextValue "E15" "  -6.24%  "

Set-TextValue "D16" "6.148"
Set-TextValue "E16" "  -6.14%  "

Set-TextValue "D17" "29.113.79"
Set-TextValue "E17" "  -1.73%  "

Set-TextValue "D18" "0.000008225"
Set-TextValue "E18" "  -1.54%  "

Set-TextValue "D19" "2.077.39"
Set-TextValue "E19" "  -1.40%  "

Set-TextValue "D20" "228.52"
Set-TextValue "E20" "  -5.51%  "

Set-TextValue "D21" "12.46"
Set-TextValue "E21" "  -2.45%  "

Set-TextValue "D22" "0.9992"
Set-TextValue "E22" "  -0.09%  "

Set-TextValue "D23" "7.330"
Set-TextValue "E23" "  -3.71%  "

Set-TextValue "D24" "1.000"
Set-TextValue "E24" "  -0.02%  "

Set-TextValue "D25" "160.94"
Set-TextValue "E25" "  +0.80%  "

Set-TextValue "D26" "0.1437"
Set-TextValue "E26" "  -4.70%  "

Set-TextValue "D27" "8.690"
Set-TextValue "E27" "  -2.48%  "

Set-TextValue "D28" "18.02"

Set-TextValue "D29" "1.500"
Set-TextValue "E29" "  -2.32%  "

Set-TextValue "D30" "4.246"

Set-TextValue "D31" "4.128"
Set-TextValue "E31" "  -1.40%  "

Set-TextValue "D32" "1.190"
Set-TextValue "E32" "  -0.31%  "

Set-TextValue "D33" "0.05380"
Set-TextValue "E33" "  +5.44%  "

Set-TextValue "D34" "0.7536"
Set-TextValue "E34" "  -3.48%  "

Set-TextValue "D35" "1.846"
Set-TextValue "E35" "  -2.93%  "

Set-TextValue "D36" "1.128"
Set-TextValue "E36" "  -2.42%  "

Set-TextValue "D37" "2.684"
Set-TextValue "E37" "  +0.01%  "

Set-TextValue "D38" "1.306.49"
Set-TextValue "E38" "  -1.43%  "

Set-TextValue "D39" "0.01814"
Set-TextValue "E39" "  -3.37%  "

Set-TextValue "D40" "2.715"
Set-TextValue "E40" "  -0.63%  "

Set-TextValue "D41" "0.9329"
Set-TextValue "E41" "  -2.93%  "

Set-TextValue "D42" "6.074"
Set-TextValue "E42" "  +3.68%  "

Set-TextValue "D43" "0.08532"
Set-TextValue "E43" "  +34.22%  "

Set-TextValue "D44" "104.95"
Set-TextValue "E44" "  -1.44%  "

Set-TextValue "D45" "0.9991"
Set-TextValue "E45" "  -0.06%  "

Set-TextValue "D46" "1.983.76"
Set-TextValue "E46" "  -1.09%  "

Set-TextValue "D47" "0.5182"
Set-TextValue "E47" "  -0.59%  "

Set-TextValue "B48" "RenderToken"
Set-TextValue "C48" "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
Set-TextValue "D48" "1.766"
Set-TextValue "E48" "  -1.15%  "

Set-TextValue "B49" "BabyDogeCoin"
Set-TextValue "C49" "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
Set-TextValue "D49" "0.00000000121"
Set-TextValue "E49" "  -4.24%  "

Set-TextValue "D50" "63.72"
Set-TextValue "E50" "  -1.22%  "

Set-TextValue "D51" "9.380"
Set-TextValue "E51" "  -4.08%  "
